$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated symbol-list values scraped on Fri Jan 20 16:51:25 UTC 2023.
# Each literal is written with a leading apostrophe so Excel stores it as
# plain text (matching the source inlineStr cells) instead of coercing
# numeric-looking strings (prices/percentages) into Number/Percent cells.

$ws.Range('D2').Value = "'291.67"
$ws.Range('E2').Value = "'0.39%"
$ws.Range('D3').Value = "'31.07"
$ws.Range('E3').Value = "'1.06%"
$ws.Range('D4').Value = "'4.951"
$ws.Range('E4').Value = "'1.64%"
$ws.Range('D5').Value = "'0.07439"
$ws.Range('E5').Value = "'2.70%"
$ws.Range('D6').Value = "'2.219"
$ws.Range('E6').Value = "'-5.47%"
$ws.Range('E7').Value = "'0.84%"
$ws.Range('B8').Value = "'MXToken"
$ws.Range('C8').Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range('D8').Value = "'0.9189"
$ws.Range('E8').Value = "'2.51%"
$ws.Range('B9').Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range('C9').Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range('D9').Value = "'0.09327"
$ws.Range('E9').Value = "'15.61%"
$ws.Range('B10').Value = "'WazirX"
$ws.Range('C10').Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range('D10').Value = "'0.1715"
$ws.Range('E10').Value = "'2.49%"
$ws.Range('B11').Value = "'MandalaExchangeToken"
$ws.Range('C11').Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range('D11').Value = "'0.08327"
$ws.Range('E11').Value = "'2.70%"
$ws.Range('B12').Value = "'BitrueCoin"
$ws.Range('C12').Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range('D12').Value = "'0.03219"
$ws.Range('E12').Value = "'4.64%"
$ws.Range('B13').Value = "'BitMartToken"
$ws.Range('C13').Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range('D13').Value = "'0.09996"
$ws.Range('E13').Value = "'-0.27%"
$ws.Range('B14').Value = "'BitForexToken"
$ws.Range('C14').Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range('D14').Value = "'0.001496"
$ws.Range('E14').Value = "'-0.58%"
$ws.Range('B15').Value = "'CoinExToken"
$ws.Range('C15').Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range('D15').Value = "'0.04499"
$ws.Range('E15').Value = "'-0.40%"
$ws.Range('D16').Value = "'0.005758"
$ws.Range('E16').Value = "'-0.58%"
$ws.Range('E17').Value = "'-0.11%"
$ws.Range('B18').Value = "'GateToken"
$ws.Range('C18').Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range('D18').Value = "'3.751"
$ws.Range('E18').Value = "'1.20%"
$ws.Range('B19').Value = "'BTSEToken"
$ws.Range('C19').Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range('D19').Value = "'2.130"
$ws.Range('E19').Value = "'2.70%"
$ws.Range('B20').Value = "'BitpandaEcosystemToken"
$ws.Range('C20').Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range('D20').Value = "'0.3330"
$ws.Range('E20').Value = "'0.39%"
$ws.Range('B21').Value = "'ProBitToken"
$ws.Range('C21').Value = "'https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range('D21').Value = "'0.1300"
$ws.Range('E21').Value = "'0.97%"
$ws.Range('B22').Value = "'MCDex"
$ws.Range('C22').Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range('D22').Value = "'4.158"
$ws.Range('E22').Value = "'4.94%"
$ws.Range('B23').Value = "'ZBToken"
$ws.Range('C23').Value = "'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range('D23').Value = "'0.2118"
$ws.Range('E23').Value = "'0.43%"
$ws.Range('D24').Value = "'0.001215"
$ws.Range('E24').Value = "'0.07%"
$ws.Range('D25').Value = "'0.004249"
$ws.Range('E25').Value = "'-3.51%"
$ws.Range('D26').Value = "'0.0001295"
$ws.Range('E26').Value = "'-0.34%"
$ws.Range('D27').Value = "'0.0003385"
$ws.Range('E27').Value = "'-0.36%"
$ws.Range('D39').Value = "'0.01597"
$ws.Range('E39').Value = "'0.71%"
$ws.Range('E40').Value = "'4.18%"
$ws.Range('D41').Value = "'0.007453"
$ws.Range('E41').Value = "'2.97%"
$ws.Range('D42').Value = "'0.009836"
$ws.Range('E42').Value = "'-1.65%"
$ws.Range('D43').Value = "'0.1353"
$ws.Range('E43').Value = "'3.08%"
$ws.Range('D44').Value = "'0.002152"
$ws.Range('E44').Value = "'3.53%"
$ws.Range('D45').Value = "'0.008766"
$ws.Range('E45').Value = "'-4.40%"
$ws.Range('D46').Value = "'0.00006109"
$ws.Range('E46').Value = "'6.77%"
$ws.Range('D47').Value = "'0.00000000748"
$ws.Range('E47').Value = "'-0.35%"
$ws.Range('D48').Value = "'2.494"
$ws.Range('E48').Value = "'11.24%"
$ws.Range('D49').Value = "'0.001994"
$ws.Range('E49').Value = "'-31.25%"
$ws.Range('D50').Value = "'0.00002094"
$ws.Range('E50').Value = "'-0.35%"
$ws.Range('D51').Value = "'0.0001995"
$ws.Range('E51').Value = "'-0.35%"

# Strip the quote-prefix style that typing an apostrophe-led value adds,
# so touched cells keep the workbook default style (no "s" attr), matching
# the original formatting exactly.
$ws.Range("B2:E51").ClearFormats()

